$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "ARR" label that was in A9, turning the row back into an empty row
$ws.Range("A9").ClearContents()

# Move the active selection to A9 (matches the saved selection state in the file)
$ws.Range("A9").Select()
